# Auto-generated: update crypto Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.885.07'
$ws.Range("E2").Value = '  -0.74%  '
$ws.Range("D3").Value = '2.970.86'
$ws.Range("E3").Value = '  -1.77%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'" + '498.71'
$ws.Range("E5").Value = '  -3.99%  '
$ws.Range("D6").Value = "'" + '137.29'
$ws.Range("E6").Value = '  -2.64%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = "'" + '0.429'
$ws.Range("E8").Value = '  -2.15%  '
$ws.Range("D9").Value = "'" + '7.33'
$ws.Range("E9").Value = '  -3.40%  '
$ws.Range("E10").Value = '  -2.56%  '
$ws.Range("D11").Value = "'" + '0.357'
$ws.Range("E11").Value = '  -1.10%  '
$ws.Range("D12").Value = '3.474.95'
$ws.Range("E12").Value = '  -1.86%  '
$ws.Range("E13").Value = '  -1.78%  '
$ws.Range("D14").Value = "'" + '25.89'
$ws.Range("E14").Value = '  -0.95%  '
$ws.Range("D15").Value = "'" + '0.0000159'
$ws.Range("E15").Value = '  -1.47%  '
$ws.Range("D16").Value = '56.938.86'
$ws.Range("E16").Value = '  -0.64%  '
$ws.Range("D17").Value = "'" + '6.06'
$ws.Range("E17").Value = '  +0.33%  '
$ws.Range("D18").Value = '2.972.50'
$ws.Range("E18").Value = '  -1.72%  '
$ws.Range("D19").Value = "'" + '12.58'
$ws.Range("E19").Value = '  -1.05%  '
$ws.Range("D20").Value = "'" + '7.82'
$ws.Range("E20").Value = '  -1.88%  '
$ws.Range("D21").Value = "'" + '319.41'
$ws.Range("E21").Value = '  -3.83%  '
$ws.Range("D22").Value = "'" + '0.999'
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("E23").Value = '  -0.75%  '
$ws.Range("D24").Value = "'" + '0.486'
$ws.Range("E24").Value = '  -0.56%  '
$ws.Range("D25").Value = "'" + '63.34'
$ws.Range("E25").Value = '  -1.11%  '
$ws.Range("D26").Value = "'" + '1.00'
$ws.Range("E26").Value = '  -0.10%  '
$ws.Range("E27").Value = '  -5.15%  '
$ws.Range("D28").Value = '0.0₃0889'
$ws.Range("E28").Value = '  -4.18%  '
$ws.Range("E29").Value = '  -3.72%  '
$ws.Range("D30").Value = "'" + '7.04'
$ws.Range("E30").Value = '  -1.61%  '
$ws.Range("E31").Value = '  -3.57%  '
$ws.Range("E32").Value = '  -6.13%  '
$ws.Range("D33").Value = "'" + '20.14'
$ws.Range("E33").Value = '  -3.26%  '
$ws.Range("D34").Value = "'" + '156.09'
$ws.Range("E34").Value = '  -1.10%  '
$ws.Range("D35").Value = "'" + '4.58'
$ws.Range("E35").Value = '  -1.55%  '
$ws.Range("D36").Value = "'" + '5.76'
$ws.Range("E36").Value = '  -0.44%  '
$ws.Range("E37").Value = '  -4.47%  '
$ws.Range("D38").Value = "'" + '24.18'
$ws.Range("E38").Value = '  -0.68%  '
$ws.Range("D39").Value = "'" + '0.0665'
$ws.Range("E39").Value = '  -2.37%  '
$ws.Range("D40").Value = '3.000.45'
$ws.Range("E40").Value = '  -1.80%  '
$ws.Range("D41").Value = "'" + '37.55'
$ws.Range("E41").Value = '  +0.44%  '
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("D43").Value = "'" + '3.73'
$ws.Range("E43").Value = '  -0.45%  '
$ws.Range("D44").Value = "'" + '0.638'
$ws.Range("E44").Value = '  -2.12%  '
$ws.Range("D45").Value = '2.199.18'
$ws.Range("E45").Value = '  -4.41%  '
$ws.Range("E46").Value = '  -3.68%  '
$ws.Range("D47").Value = "'" + '0.944'
$ws.Range("E47").Value = '  -6.87%  '
$ws.Range("D48").Value = "'" + '5.93'
$ws.Range("E48").Value = '  +0.29%  '
$ws.Range("D49").Value = "'" + '0.0234'
$ws.Range("E49").Value = '  -3.78%  '
$ws.Range("D50").Value = "'" + '19.19'
$ws.Range("E50").Value = '  -1.71%  '
$ws.Range("E51").Value = '  -11.06%  '
